$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.1239850077872064, 0.9587971593447111, 2.504666468680252, 1.582613809076697, 1.593449111284933, 51),
    @(-0.07634650817340942, 0.9644256602613136, 2.465022428948611, 1.570038989626885, 1.584102668685676, 50),
    @(0.1308097497450276, 0.9583096004403104, 2.464810574400857, 1.569971520251516, 1.580725504613128, 49),
    @(-0.02424481833496826, 0.9746006404459205, 2.576612619674879, 1.605183048650489, 1.621984539672506, 48),
    @(0.1313683227808764, 1.001393086257909, 2.609607854279241, 1.615428071527557, 1.627484456937081, 47),
    @(0.003180273963439154, 0.9469229433371078, 2.518249158100375, 1.586899227455977, 1.604431335894871, 46),
    @(0.1362810442212024, 0.9540930716873162, 2.589108079185797, 1.60907056376835, 1.621405815584749, 45),
    @(-0.04071899536153733, 0.8760439952630932, 2.432158382960957, 1.559537874808097, 1.577029996622433, 44),
    @(0.07889961332241763, 0.9613424106959066, 2.595247230407774, 1.610977104246915, 1.628086442519692, 43),
    @(-0.01422600529044099, 0.8653212001480697, 2.498623644143271, 1.580703528225097, 1.599799477728101, 42)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
